$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date from 2023-10-22 to 2023-10-25 for rows 2-18
$newDate = Get-Date -Year 2023 -Month 10 -Day 25 -Hour 0 -Minute 0 -Second 0

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
